$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" note with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.12 = 24403.67 pesos`n✅ 24403.67 pesos = 6.09 = 955.93 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas sheet: update the N10/O10/N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 163.5
$ws2.Range("O10").Value = 3990
$ws2.Range("N12").Value = 4008.5
$ws2.Range("O12").Value = 157.02
